$wb = $excel.ActiveWorkbook

# --- Rename "Sheet6" to "Demo" ---
$demo = $wb.Worksheets.Item("Sheet6")
$demo.Name = "Demo"

# --- Hide "Sheet1" ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Visible = $false

# --- Update the Demo sheet's selection / styling / new unit string ---
$demo.Activate()

# Align street name + house number cells to the right (matching the other value cells)
$demo.Range("E5").HorizontalAlignment = -4152   # xlRight
$demo.Range("E6").HorizontalAlignment = -4152   # xlRight

# Add the "m2" unit label next to the surface area value
$demo.Range("F10").Value = "m2"

# Move the active selection to H13
$demo.Range("H13").Select()

# First visible tab should be the (now hidden) Sheet1's neighbour; scroll tab strip so Demo is first shown
$wb.Windows.Item(1).ScrollWorkbookTabs(1)
